$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give every data row (beyond the first, which already has a unique id) its
# own unique user_id / transaction_id, instead of sharing the same two
# strings across all 103 remaining rows. This is the fix described in the
# commit message: "unique id for each row so that they will work with
# tests" (avoids race conditions on rollback when rows collide on id).
for ($r = 3; $r -le 105; $r++) {
    $userIdSuffix = 7621 + ($r - 3)
    $transactionIdSuffix = 2301 + ($r - 3)
    $ws.Cells.Item($r, 3).Value = "988-90-$userIdSuffix"
    $ws.Cells.Item($r, 4).Value = "293-64-$transactionIdSuffix"
}

# Column widths grew (auto-fit) because the new id strings changed the
# content used to size the columns. Apply the resulting widths directly
# (Excel's ColumnWidth property is ~5/6 of a character narrower than the
# width value persisted in the sheet XML).
$ws.Columns.Item(1).ColumnWidth = 22.8571428571429 - 0.8333333333333334
$ws.Columns.Item(2).ColumnWidth = 19.4387755102041 - 0.8333333333333334
$ws.Columns.Item(3).ColumnWidth = 15.3010204081633 - 0.8333333333333334
$ws.Columns.Item(4).ColumnWidth = 16.3775510204082 - 0.8333333333333334
$ws.Columns.Item(5).ColumnWidth = 20.3367346938776 - 0.8333333333333334

# The saved file's view ends up scrolled down with E104 (the last edited
# cell) selected.
$ws.Range("E104").Select()
